$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-11 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-12 Monday", 2) | Out-Null
$d.Content.Find.Execute("56×88=4928", $true, $false, $false, $false, $false, $true, 1, $false, "79×81=6399", 2) | Out-Null
$d.Content.Find.Execute("35×13=455", $true, $false, $false, $false, $false, $true, 1, $false, "89×40=3560", 2) | Out-Null
$d.Content.Find.Execute("16×28=448", $true, $false, $false, $false, $false, $true, 1, $false, "92×14=1288", 2) | Out-Null
$d.Content.Find.Execute("54×92=4968", $true, $false, $false, $false, $false, $true, 1, $false, "40×85=3400", 2) | Out-Null
$d.Content.Find.Execute("89×16=1424", $true, $false, $false, $false, $false, $true, 1, $false, "32×18=576", 2) | Out-Null
$d.Content.Find.Execute("60×98=5880", $true, $false, $false, $false, $false, $true, 1, $false, "70×48=3360", 2) | Out-Null
$d.Content.Find.Execute("21×49=1029", $true, $false, $false, $false, $false, $true, 1, $false, "92×91=8372", 2) | Out-Null
$d.Content.Find.Execute("99×65=6435", $true, $false, $false, $false, $false, $true, 1, $false, "56×50=2800", 2) | Out-Null
$d.Content.Find.Execute("27×40=1080", $true, $false, $false, $false, $false, $true, 1, $false, "30×48=1440", 2) | Out-Null
$d.Content.Find.Execute("99×19=1881", $true, $false, $false, $false, $false, $true, 1, $false, "69×95=6555", 2) | Out-Null
$d.Content.Find.Execute("17×79=1343", $true, $false, $false, $false, $false, $true, 1, $false, "13×89=1157", 2) | Out-Null
$d.Content.Find.Execute("89×11=979", $true, $false, $false, $false, $false, $true, 1, $false, "82×49=4018", 2) | Out-Null
$d.Content.Find.Execute("28×58=1624", $true, $false, $false, $false, $false, $true, 1, $false, "37×72=2664", 2) | Out-Null
$d.Content.Find.Execute("30×57=1710", $true, $false, $false, $false, $false, $true, 1, $false, "26×46=1196", 2) | Out-Null
$d.Content.Find.Execute("70×28=1960", $true, $false, $false, $false, $false, $true, 1, $false, "93×24=2232", 2) | Out-Null
$d.Content.Find.Execute("17×28=476", $true, $false, $false, $false, $false, $true, 1, $false, "90×44=3960", 2) | Out-Null
$d.Content.Find.Execute("94×39=3666", $true, $false, $false, $false, $false, $true, 1, $false, "75×89=6675", 2) | Out-Null
$d.Content.Find.Execute("50×38=1900", $true, $false, $false, $false, $false, $true, 1, $false, "74×47=3478", 2) | Out-Null
$d.Content.Find.Execute("51×64=3264", $true, $false, $false, $false, $false, $true, 1, $false, "39×68=2652", 2) | Out-Null
$d.Content.Find.Execute("46×69=3174", $true, $false, $false, $false, $false, $true, 1, $false, "51×24=1224", 2) | Out-Null
$d.Content.Find.Execute("68×92=6256", $true, $false, $false, $false, $false, $true, 1, $false, "30×76=2280", 2) | Out-Null
$d.Content.Find.Execute("29×21=609", $true, $false, $false, $false, $false, $true, 1, $false, "34×88=2992", 2) | Out-Null
$d.Content.Find.Execute("98×26=2548", $true, $false, $false, $false, $false, $true, 1, $false, "45×63=2835", 2) | Out-Null
$d.Content.Find.Execute("23×78=1794", $true, $false, $false, $false, $false, $true, 1, $false, "49×77=3773", 2) | Out-Null
$d.Content.Find.Execute("85×77=6545", $true, $false, $false, $false, $false, $true, 1, $false, "31×23=713", 2) | Out-Null

Write-Host "Done replacing text."
